$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection to G11
$ws.Range("G11").Select()

# Update G2 value
$ws.Range("G2").Value = 0.75600000000000001

# Update G5 value
$ws.Range("G5").Value = 0.29099999999999998

# Update G7: apply the same number format as the rest of column G (percentage, 2 decimals)
# and set its value to 1.2999999999999999E-2 (1.3%)
$ws.Range("G7").NumberFormat = "0.00%"
$ws.Range("G7").Value = 0.012999999999999999

# Update G9 value
$ws.Range("G9").Value = 0.748
